# "automation cho chuc nang update"
# Adds a second worksheet ("Sheet2") right after "Sheet1", reusing Sheet1's
# formatting (column style / fonts) via Copy(), then replacing the copied
# content with the new test data, and updates which sheet/cell is selected
# on each sheet to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1 ends up with no tabSelected flag and its selection moves to A2 -
# do this before creating Sheet2 so Sheet2 ends up as the final active tab.
$ws1.Range("A2").Select() | Out-Null

# Duplicate Sheet1 right after itself - this carries over the same column
# formatting/styles (14pt "Times New Roman", style index 1) that the new
# sheet ends up using for every cell it holds.
$ws1.Copy($null, $ws1) | Out-Null

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Drop the hyperlink + all inherited values from Sheet1, then trim the
# extra copied rows (7:8) that Sheet2 doesn't need.
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Range("A2:D8").Clear() | Out-Null
$ws2.Rows("7:8").Delete() | Out-Null

# New Sheet2 content.
$ws2.Range("A2").Value = "invalid"
$ws2.Range("B3").Value = "xem phim"
$ws2.Range("B5").Value = "xem"
$ws2.Range("C3").Value = "hai"
$ws2.Range("C4").Value = "hai"
$ws2.Range("D6").Value = "invalid"

# Sheet2 is the active sheet/tab, with B5 selected.
$ws2.Range("B5").Select() | Out-Null
